$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2078.7847
$ws.Range("J17").Value = 2078.7847
$ws.Range("L17").Value = 6236.3541
$ws.Range("N17").Value = -6572.3541
$ws.Range("H40").Value = 7422.4
$ws.Range("J40").Value = 10657.5
$ws.Range("L40").Value = 10657.5
$ws.Range("N40").Value = -11007.5
$ws.Range("H58").Value = 6350.25
$ws.Range("I58").Value = 365.83334
$ws.Range("K58").Value = 1097.50002
$ws.Range("M58").Value = -947.5000199999999
$ws.Range("H86").Value = 2927303
$ws.Range("J86").Value = 5267257.5
$ws.Range("L86").Value = 5267257.5
$ws.Range("N86").Value = -5269503.5
$ws.Range("H89").Value = 2927303
$ws.Range("J89").Value = 5267257.5
$ws.Range("L89").Value = 26336287.5
$ws.Range("N89").Value = -26347519.5
$ws.Range("H106").Value = 3515.2104
$ws.Range("I106").Value = 3433
$ws.Range("K106").Value = 3433
$ws.Range("M106").Value = -2802
$ws.Range("H137").Value = 1776.1666
$ws.Range("J137").Value = 2621.3333
$ws.Range("L137").Value = 7863.999899999999
$ws.Range("N137").Value = -12963.9999
$ws.Range("H138").Value = 6543.4746
$ws.Range("I138").Value = 3449.0625
$ws.Range("J138").Value = 7694.884
$ws.Range("K138").Value = 10347.1875
$ws.Range("L138").Value = 23084.652
$ws.Range("M138").Value = -5207.1875
$ws.Range("N138").Value = -33364.652

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7671.773
$ws.Range("I32").Value = 7599.0244
$ws.Range("K32").Value = 7599.0244
$ws.Range("M32").Value = -7312.0244
$ws.Range("H45").Value = 3616.4546
$ws.Range("J45").Value = 5999.75
$ws.Range("L45").Value = 5999.75
$ws.Range("N45").Value = -6753.75
$ws.Range("H102").Value = 2007.8462
$ws.Range("I102").Value = 2081.4
$ws.Range("J102").Value = 1762.6666
$ws.Range("K102").Value = 2081.4
$ws.Range("L102").Value = 1762.6666
$ws.Range("M102").Value = -459.4000000000001
$ws.Range("N102").Value = -5006.6666
$ws.Range("H132").Value = 2965.3655
$ws.Range("I132").Value = 2766.6445
$ws.Range("K132").Value = 8299.933499999999
$ws.Range("M132").Value = -5769.933499999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492
$ws.Range("H105").Value = 73954.5
$ws.Range("J105").Value = 1538.75
$ws.Range("L105").Value = 1538.75
$ws.Range("N105").Value = -5032.75
$ws.Range("H134").Value = 23791.521
$ws.Range("I134").Value = 3056.0513
$ws.Range("K134").Value = 9168.153900000001
$ws.Range("M134").Value = -6633.153900000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 46066.695
$ws.Range("I31").Value = 1533.4615
$ws.Range("J31").Value = 103959.9
$ws.Range("K31").Value = 1533.4615
$ws.Range("L31").Value = 103959.9
$ws.Range("M31").Value = -1238.4615
$ws.Range("N31").Value = -104549.9
$ws.Range("H34").Value = 46066.695
$ws.Range("I34").Value = 1533.4615
$ws.Range("J34").Value = 103959.9
$ws.Range("K34").Value = 1533.4615
$ws.Range("L34").Value = 103959.9
$ws.Range("M34").Value = -1331.4615
$ws.Range("N34").Value = -104363.9
$ws.Range("H99").Value = 7549.4443
$ws.Range("I99").Value = 8739
$ws.Range("K99").Value = 8739
$ws.Range("M99").Value = -7241
$ws.Range("H126").Value = 7549.4443
$ws.Range("I126").Value = 8739
$ws.Range("K126").Value = 26217
$ws.Range("M126").Value = -23747
$ws.Range("H132").Value = 1688.1177
$ws.Range("I132").Value = 1370.75
$ws.Range("J132").Value = 3169.1667
$ws.Range("K132").Value = 4112.25
$ws.Range("L132").Value = 9507.500100000001
$ws.Range("M132").Value = -1582.25
$ws.Range("N132").Value = -14567.5001
$ws.Range("H134").Value = 252088.62
$ws.Range("I134").Value = 2090.3713
$ws.Range("J134").Value = 2002076.4
$ws.Range("K134").Value = 6271.113899999999
$ws.Range("L134").Value = 6006229.199999999
$ws.Range("M134").Value = -3736.113899999999
$ws.Range("N134").Value = -6011299.199999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H118").Value = 8136.242
$ws.Range("I118").Value = 1500
$ws.Range("J118").Value = 8799.866
$ws.Range("K118").Value = 4500
$ws.Range("L118").Value = 26399.598
$ws.Range("M118").Value = -3257
$ws.Range("N118").Value = -28885.598
$ws.Range("H119").Value = 21066
$ws.Range("I119").Value = 10507
$ws.Range("J119").Value = 31625
$ws.Range("K119").Value = 31521
$ws.Range("L119").Value = 94875
$ws.Range("M119").Value = -26683
$ws.Range("N119").Value = -104551
$ws.Range("H132").Value = 461408.4
$ws.Range("I132").Value = 144614.14
$ws.Range("J132").Value = 591853.1
$ws.Range("K132").Value = 1301527.26
$ws.Range("L132").Value = 5326677.899999999
$ws.Range("M132").Value = -1298997.26
$ws.Range("N132").Value = -5331737.899999999
$ws.Range("H140").Value = 5161.433
$ws.Range("I140").Value = 3380.6667
$ws.Range("J140").Value = 6942.2
$ws.Range("K140").Value = 10142.0001
$ws.Range("L140").Value = 20826.6
$ws.Range("M140").Value = -4962.000100000001
$ws.Range("N140").Value = -31186.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2173.2173
$ws.Range("I102").Value = 1314.8182
$ws.Range("K102").Value = 1314.8182
$ws.Range("M102").Value = 307.1818000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H7").Value = 4799.5356
$ws.Range("I7").Value = 4404.421
$ws.Range("K7").Value = 4404.421
$ws.Range("M7").Value = -4292.421
$ws.Range("H93").Value = 3232.7036
$ws.Range("I93").Value = 2842.5454
$ws.Range("K93").Value = 2842.5454
$ws.Range("M93").Value = -1594.5454
$ws.Range("H126").Value = 4799.5356
$ws.Range("I126").Value = 4404.421
$ws.Range("K126").Value = 13213.263
$ws.Range("M126").Value = -10743.263
$ws.Range("H132").Value = 3181.5454
$ws.Range("I132").Value = 2999.8
$ws.Range("K132").Value = 8999.400000000001
$ws.Range("M132").Value = -6469.400000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1557.4546
$ws.Range("I107").Value = 1729.1666
$ws.Range("J107").Value = 784.75
$ws.Range("K107").Value = 5187.4998
$ws.Range("L107").Value = 2354.25
$ws.Range("M107").Value = -3267.4998
$ws.Range("N107").Value = -6194.25
$ws.Range("H122").Value = 43482228
$ws.Range("I122").Value = 66669732
$ws.Range("K122").Value = 200009196
$ws.Range("M122").Value = -200006746
$ws.Range("H132").Value = 15371.925
$ws.Range("I132").Value = 2253.1384
$ws.Range("K132").Value = 6759.415199999999
$ws.Range("M132").Value = -4229.415199999999
$ws.Range("H136").Value = 53863.332
$ws.Range("I136").Value = 2552
$ws.Range("J136").Value = 502837.5
$ws.Range("K136").Value = 7656
$ws.Range("L136").Value = 1508512.5
$ws.Range("M136").Value = -5106
$ws.Range("N136").Value = -1513612.5
